$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")

# Rename existing "accountOfficer" row to "accountOfficer1"
$ws.Cells.Item(8, 1).Value = "accountOfficer1"

# Add a new row for the second account officer / examiner of accounts workflow step
$ws.Cells.Item(9, 1).Value = "accountOfficer2"
$ws.Cells.Item(9, 2).Value = "ACCOUNTS"
$ws.Cells.Item(9, 3).Value = "Examiner of Accounts"
$ws.Cells.Item(9, 4).Value = "D Ramachandra Reddy ~ ACC_EOA_1"

# Make approvalDetails the active sheet / tab, with A5 selected
$ws.Activate()
$ws.Range("A5").Select()
